$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1465.24
$ws.Cells.Item(15, 9).Value = 1465.24
$ws.Cells.Item(15, 11).Value = 4395.72
$ws.Cells.Item(15, 13).Value = -4226.72
$ws.Cells.Item(40, 8).Value = 4975.05
$ws.Cells.Item(40, 9).Value = 2346.2307
$ws.Cells.Item(40, 11).Value = 2346.2307
$ws.Cells.Item(40, 13).Value = -2171.2307
$ws.Cells.Item(55, 8).Value = 56447
$ws.Cells.Item(55, 9).Value = 206.66667
$ws.Cells.Item(55, 10).Value = 67695.07000000001
$ws.Cells.Item(55, 11).Value = 206.66667
$ws.Cells.Item(55, 12).Value = 67695.07000000001
$ws.Cells.Item(55, 13).Value = 7.333329999999989
$ws.Cells.Item(55, 14).Value = -68123.07000000001
$ws.Cells.Item(100, 8).Value = 3214.5217
$ws.Cells.Item(100, 9).Value = 2366.3845
$ws.Cells.Item(100, 10).Value = 4317.1
$ws.Cells.Item(100, 11).Value = 2366.3845
$ws.Cells.Item(100, 12).Value = 4317.1
$ws.Cells.Item(100, 13).Value = -1825.3845
$ws.Cells.Item(100, 14).Value = -5399.1
$ws.Cells.Item(103, 8).Value = 0
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 11).Value = 0
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(112, 8).Value = 3240.7693
$ws.Cells.Item(112, 9).Value = 1333
$ws.Cells.Item(112, 10).Value = 3399.75
$ws.Cells.Item(112, 11).Value = 3999
$ws.Cells.Item(112, 12).Value = 10199.25
$ws.Cells.Item(112, 13).Value = -2891
$ws.Cells.Item(112, 14).Value = -12415.25
$ws.Cells.Item(132, 8).Value = 3478.4194
$ws.Cells.Item(132, 9).Value = 3606.423
$ws.Cells.Item(132, 11).Value = 10819.269
$ws.Cells.Item(132, 13).Value = -8289.269
$ws.Cells.Item(137, 8).Value = 6020.227
$ws.Cells.Item(137, 10).Value = 6373
$ws.Cells.Item(137, 12).Value = 19119
$ws.Cells.Item(137, 14).Value = -24219
$ws.Cells.Item(138, 8).Value = 8343.75
$ws.Cells.Item(138, 9).Value = 6101.75
$ws.Cells.Item(138, 10).Value = 9688.950000000001
$ws.Cells.Item(138, 11).Value = 18305.25
$ws.Cells.Item(138, 12).Value = 29066.85
$ws.Cells.Item(138, 13).Value = -13165.25
$ws.Cells.Item(138, 14).Value = -39346.85000000001
$ws.Cells.Item(103, 13).ClearContents()
$ws.Cells.Item(103, 14).ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2377
$ws.Cells.Item(2, 9).Value = 1755.8
$ws.Cells.Item(2, 11).Value = 1755.8
$ws.Cells.Item(2, 13).Value = -1642.8
$ws.Cells.Item(32, 8).Value = 2990.5469
$ws.Cells.Item(32, 9).Value = 2317.9673
$ws.Cells.Item(32, 11).Value = 2317.9673
$ws.Cells.Item(32, 13).Value = -2030.9673
$ws.Cells.Item(45, 8).Value = 7166.222
$ws.Cells.Item(45, 9).Value = 3499.4285
$ws.Cells.Item(45, 10).Value = 20000
$ws.Cells.Item(45, 11).Value = 3499.4285
$ws.Cells.Item(45, 12).Value = 20000
$ws.Cells.Item(45, 13).Value = -3122.4285
$ws.Cells.Item(45, 14).Value = -20754
$ws.Cells.Item(116, 8).Value = 2377
$ws.Cells.Item(116, 9).Value = 1755.8
$ws.Cells.Item(116, 11).Value = 1755.8
$ws.Cells.Item(116, 13).Value = 538.2
$ws.Cells.Item(132, 8).Value = 6776.475
$ws.Cells.Item(132, 9).Value = 2263.08
$ws.Cells.Item(132, 10).Value = 14298.8
$ws.Cells.Item(132, 11).Value = 6789.24
$ws.Cells.Item(132, 12).Value = 42896.39999999999
$ws.Cells.Item(132, 13).Value = -4259.24
$ws.Cells.Item(132, 14).Value = -47956.39999999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2377
$ws.Cells.Item(3, 9).Value = 1755.8
$ws.Cells.Item(3, 11).Value = 1755.8
$ws.Cells.Item(3, 13).Value = -1641.8
$ws.Cells.Item(105, 8).Value = 15154480
$ws.Cells.Item(105, 9).Value = 2522
$ws.Cells.Item(105, 10).Value = 27781112
$ws.Cells.Item(105, 11).Value = 2522
$ws.Cells.Item(105, 12).Value = 27781112
$ws.Cells.Item(105, 13).Value = -775
$ws.Cells.Item(105, 14).Value = -27784606
$ws.Cells.Item(134, 8).Value = 30467.13
$ws.Cells.Item(134, 9).Value = 3420.1667
$ws.Cells.Item(134, 10).Value = 76833.36
$ws.Cells.Item(134, 11).Value = 10260.5001
$ws.Cells.Item(134, 12).Value = 230500.08
$ws.Cells.Item(134, 13).Value = -7725.500100000001
$ws.Cells.Item(134, 14).Value = -235570.08

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4781.353
$ws.Cells.Item(31, 9).Value = 2557.2
$ws.Cells.Item(31, 10).Value = 5708.0835
$ws.Cells.Item(31, 11).Value = 2557.2
$ws.Cells.Item(31, 12).Value = 5708.0835
$ws.Cells.Item(31, 13).Value = -2262.2
$ws.Cells.Item(31, 14).Value = -6298.0835
$ws.Cells.Item(34, 8).Value = 4781.353
$ws.Cells.Item(34, 9).Value = 2557.2
$ws.Cells.Item(34, 10).Value = 5708.0835
$ws.Cells.Item(34, 11).Value = 2557.2
$ws.Cells.Item(34, 12).Value = 5708.0835
$ws.Cells.Item(34, 13).Value = -2355.2
$ws.Cells.Item(34, 14).Value = -6112.0835
$ws.Cells.Item(58, 8).Value = 281168.1
$ws.Cells.Item(58, 9).Value = 418481.5
$ws.Cells.Item(58, 10).Value = 6541.25
$ws.Cells.Item(58, 11).Value = 418481.5
$ws.Cells.Item(58, 12).Value = 6541.25
$ws.Cells.Item(58, 13).Value = -418278.5
$ws.Cells.Item(58, 14).Value = -6947.25
$ws.Cells.Item(99, 8).Value = 3775.5667
$ws.Cells.Item(99, 10).Value = 5000.8335
$ws.Cells.Item(99, 12).Value = 5000.8335
$ws.Cells.Item(99, 14).Value = -7996.8335
$ws.Cells.Item(126, 8).Value = 3775.5667
$ws.Cells.Item(126, 10).Value = 5000.8335
$ws.Cells.Item(126, 12).Value = 15002.5005
$ws.Cells.Item(126, 14).Value = -19942.5005
$ws.Cells.Item(132, 8).Value = 4306.7827
$ws.Cells.Item(132, 9).Value = 4468.8203
$ws.Cells.Item(132, 10).Value = 3404
$ws.Cells.Item(132, 11).Value = 13406.4609
$ws.Cells.Item(132, 12).Value = 10212
$ws.Cells.Item(132, 13).Value = -10876.4609
$ws.Cells.Item(132, 14).Value = -15272
$ws.Cells.Item(134, 8).Value = 480002.47
$ws.Cells.Item(134, 9).Value = 3203.4666
$ws.Cells.Item(134, 10).Value = 1672000
$ws.Cells.Item(134, 11).Value = 9610.399800000001
$ws.Cells.Item(134, 12).Value = 5016000
$ws.Cells.Item(134, 13).Value = -7075.399800000001
$ws.Cells.Item(134, 14).Value = -5021070
$ws.Cells.Item(136, 8).Value = 281168.1
$ws.Cells.Item(136, 9).Value = 418481.5
$ws.Cells.Item(136, 10).Value = 6541.25
$ws.Cells.Item(136, 11).Value = 1255444.5
$ws.Cells.Item(136, 12).Value = 19623.75
$ws.Cells.Item(136, 13).Value = -1252894.5
$ws.Cells.Item(136, 14).Value = -24723.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 68356.164
$ws.Cells.Item(5, 9).Value = 161420
$ws.Cells.Item(5, 10).Value = 1882
$ws.Cells.Item(5, 11).Value = 484260
$ws.Cells.Item(5, 12).Value = 5646
$ws.Cells.Item(5, 13).Value = -484148
$ws.Cells.Item(5, 14).Value = -5870
$ws.Cells.Item(56, 8).Value = 6905.625
$ws.Cells.Item(56, 9).Value = 6905.625
$ws.Cells.Item(56, 11).Value = 6905.625
$ws.Cells.Item(56, 13).Value = -6375.625
$ws.Cells.Item(68, 8).Value = 2605.6667
$ws.Cells.Item(68, 9).Value = 2327.4666
$ws.Cells.Item(68, 10).Value = 3996.6667
$ws.Cells.Item(68, 11).Value = 6982.399800000001
$ws.Cells.Item(68, 12).Value = 11990.0001
$ws.Cells.Item(68, 13).Value = -6171.399800000001
$ws.Cells.Item(68, 14).Value = -13612.0001
$ws.Cells.Item(71, 8).Value = 2605.6667
$ws.Cells.Item(71, 9).Value = 2327.4666
$ws.Cells.Item(71, 10).Value = 3996.6667
$ws.Cells.Item(71, 11).Value = 20947.1994
$ws.Cells.Item(71, 12).Value = 35970.0003
$ws.Cells.Item(71, 13).Value = -16891.1994
$ws.Cells.Item(71, 14).Value = -44082.0003
$ws.Cells.Item(131, 8).Value = 8556.857
$ws.Cells.Item(131, 9).Value = 2199
$ws.Cells.Item(131, 10).Value = 9616.5
$ws.Cells.Item(131, 11).Value = 6597
$ws.Cells.Item(131, 12).Value = 28849.5
$ws.Cells.Item(131, 13).Value = -1557
$ws.Cells.Item(131, 14).Value = -38929.5
$ws.Cells.Item(135, 8).Value = 68356.164
$ws.Cells.Item(135, 9).Value = 161420
$ws.Cells.Item(135, 10).Value = 1882
$ws.Cells.Item(135, 11).Value = 1452780
$ws.Cells.Item(135, 12).Value = 16938
$ws.Cells.Item(135, 13).Value = -1450245
$ws.Cells.Item(135, 14).Value = -22008
$ws.Cells.Item(137, 8).Value = 3674.353
$ws.Cells.Item(137, 10).Value = 4998.2856
$ws.Cells.Item(137, 12).Value = 14994.8568
$ws.Cells.Item(137, 14).Value = -25194.8568

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 958.2963
$ws.Cells.Item(97, 9).Value = 1113.2
$ws.Cells.Item(97, 10).Value = 515.7143
$ws.Cells.Item(97, 11).Value = 1113.2
$ws.Cells.Item(97, 12).Value = 515.7143
$ws.Cells.Item(97, 13).Value = -617.2
$ws.Cells.Item(97, 14).Value = -1507.7143
$ws.Cells.Item(99, 8).Value = 2743
$ws.Cells.Item(99, 9).Value = 2743
$ws.Cells.Item(99, 11).Value = 2743
$ws.Cells.Item(99, 13).Value = -497
$ws.Cells.Item(132, 8).Value = 34001
$ws.Cells.Item(132, 9).Value = 2677.2693
$ws.Cells.Item(132, 10).Value = 169737.17
$ws.Cells.Item(132, 11).Value = 8031.8079
$ws.Cells.Item(132, 12).Value = 509211.51
$ws.Cells.Item(132, 13).Value = -5501.8079
$ws.Cells.Item(132, 14).Value = -514271.51

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 522.7037
$ws.Cells.Item(55, 9).Value = 326.94736
$ws.Cells.Item(55, 11).Value = 326.94736
$ws.Cells.Item(55, 13).Value = -153.94736
$ws.Cells.Item(61, 8).Value = 7853.2
$ws.Cells.Item(61, 9).Value = 6422
$ws.Cells.Item(61, 10).Value = 10000
$ws.Cells.Item(61, 11).Value = 6422
$ws.Cells.Item(61, 12).Value = 10000
$ws.Cells.Item(61, 13).Value = -6220
$ws.Cells.Item(61, 14).Value = -10404
$ws.Cells.Item(113, 8).Value = 7853.2
$ws.Cells.Item(113, 9).Value = 6422
$ws.Cells.Item(113, 10).Value = 10000
$ws.Cells.Item(113, 11).Value = 6422
$ws.Cells.Item(113, 12).Value = 10000
$ws.Cells.Item(113, 13).Value = -4252
$ws.Cells.Item(113, 14).Value = -14340

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2914.6316
$ws.Cells.Item(81, 9).Value = 2914.6316
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 5829.2632
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 13).Value = -4768.2632
$ws.Cells.Item(84, 8).Value = 2914.6316
$ws.Cells.Item(84, 9).Value = 2914.6316
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 29146.316
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = -23842.316
$ws.Cells.Item(136, 8).Value = 637387.5600000001
$ws.Cells.Item(136, 10).Value = 207675.3
$ws.Cells.Item(136, 12).Value = 623025.8999999999
$ws.Cells.Item(136, 14).Value = -628125.8999999999
$ws.Cells.Item(81, 14).ClearContents()
$ws.Cells.Item(84, 14).ClearContents()
